$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 36.75793933333333
$ws.Range("H2").Value = 110.273818
$ws.Range("I2").Value = 0.9858943139827973
$ws.Range("J2").Value = 0.9858943139827971
$ws.Range("M2").Value = 3.241087666666667
$ws.Range("N2").Value = 9.723262999999999
$ws.Range("O2").Value = 0.02486257877280725
$ws.Range("P2").Value = 0.02486257877280725
$ws.Range("Q2").Value = 119.1357038253482
$ws.Range("R2").Value = 1072.221334428134
$ws.Range("S2").Value = 0.02451187504306006
$ws.Range("T2").Value = 0.02451187504306006
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 36.75793933333333
$ws.Range("H3").Value = 110.273818
$ws.Range("I3").Value = 0.9858943139827973
$ws.Range("J3").Value = 0.9858943139827971
$ws.Range("O3").Value = 0.02096124117795788
$ws.Range("P3").Value = 0.02096124117795788
$ws.Range("Q3").Value = 100.4414000497873
$ws.Range("R3").Value = 903.9726004480859
$ws.Range("S3").Value = 0.02066556849137075
$ws.Range("T3").Value = 0.02066556849137074
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 36.75793933333333
$ws.Range("H4").Value = 110.273818
$ws.Range("I4").Value = 0.9858943139827973
$ws.Range("J4").Value = 0.9858943139827971
$ws.Range("M4").Value = 124.3864796666667
$ws.Range("N4").Value = 373.159439
$ws.Range("O4").Value = 0.9541761800492348
$ws.Range("P4").Value = 0.9541761800492349
$ws.Range("Q4").Value = 4572.190673474233
$ws.Range("R4").Value = 41149.7160612681
$ws.Range("S4").Value = 0.9407168704483664
$ws.Range("T4").Value = 0.9407168704483664
$ws.Range("I5").Value = 0.001251989679428792
$ws.Range("J5").Value = 0.001251989679428792
$ws.Range("M5").Value = 3.241087666666667
$ws.Range("N5").Value = 9.723262999999999
$ws.Range("O5").Value = 0.02486257877280725
$ws.Range("P5").Value = 0.02486257877280725
$ws.Range("Q5").Value = 0.1512907311923333
$ws.Range("R5").Value = 1.361616580731
$ws.Range("S5").Value = [double]"3.112769202754004E-05"
$ws.Range("T5").Value = [double]"3.112769202754004E-05"
$ws.Range("I6").Value = 0.001251989679428792
$ws.Range("J6").Value = 0.001251989679428792
$ws.Range("O6").Value = 0.02096124117795788
$ws.Range("P6").Value = 0.02096124117795788
$ws.Range("S6").Value = [double]"2.624325762282109E-05"
$ws.Range("T6").Value = [double]"2.624325762282108E-05"
$ws.Range("I7").Value = 0.001251989679428792
$ws.Range("J7").Value = 0.001251989679428792
$ws.Range("M7").Value = 124.3864796666667
$ws.Range("N7").Value = 373.159439
$ws.Range("O7").Value = 0.9541761800492348
$ws.Range("P7").Value = 0.9541761800492349
$ws.Range("Q7").Value = 5.806236484360333
$ws.Range("R7").Value = 52.256128359243
$ws.Range("S7").Value = 0.001194618729778431
$ws.Range("T7").Value = 0.001194618729778431
$ws.Range("G8").Value = 0.4792353333333333
$ws.Range("H8").Value = 1.437706
$ws.Range("I8").Value = 0.01285369633777395
$ws.Range("J8").Value = 0.01285369633777395
$ws.Range("M8").Value = 3.241087666666667
$ws.Range("N8").Value = 9.723262999999999
$ws.Range("O8").Value = 0.02486257877280725
$ws.Range("P8").Value = 0.02486257877280725
$ws.Range("Q8").Value = 1.553243728297555
$ws.Range("R8").Value = 13.979193554678
$ws.Range("S8").Value = 0.0003195760377196489
$ws.Range("T8").Value = 0.0003195760377196489
$ws.Range("G9").Value = 0.4792353333333333
$ws.Range("H9").Value = 1.437706
$ws.Range("I9").Value = 0.01285369633777395
$ws.Range("J9").Value = 0.01285369633777395
$ws.Range("O9").Value = 0.02096124117795788
$ws.Range("P9").Value = 0.02096124117795788
$ws.Range("Q9").Value = 1.309514861451333
$ws.Range("R9").Value = 11.785633753062
$ws.Range("S9").Value = 0.0002694294289643138
$ws.Range("T9").Value = 0.0002694294289643138
$ws.Range("G10").Value = 0.4792353333333333
$ws.Range("H10").Value = 1.437706
$ws.Range("I10").Value = 0.01285369633777395
$ws.Range("J10").Value = 0.01285369633777395
$ws.Range("M10").Value = 124.3864796666667
$ws.Range("N10").Value = 373.159439
$ws.Range("O10").Value = 0.9541761800492348
$ws.Range("P10").Value = 0.9541761800492349
$ws.Range("Q10").Value = 59.61039604521488
$ws.Range("R10").Value = 536.493564406934
$ws.Range("S10").Value = 0.01226469087108999
$ws.Range("T10").Value = 0.01226469087108999
